$newText = @"
<Bold>e038 Orders Phase</Bold> 
<InlineUIContainer><Button Content='r4.73' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>    
<LineBreak/><LineBreak/>
Click on blue squares to open hatches. Click on Open Hatch marker to close.
 <LineBreak/><LineBreak/>
Click crew member action boxes to select from a pull down to assign crew actions for each crew member per 
<InlineUIContainer><Button Content='r8.0' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>.
<LineBreak/><LineBreak/>
Click on the appropriate Gun Load box to set the Gun Reload marker per 
<InlineUIContainer><Button Content='r5.23' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>.
<LineBreak/><LineBreak/>
Click the Gun Reload marker/button if you want the reload to come from the ready rack. A Ready Rack Ammo Reload marker is added per  
<InlineUIContainer><Button Content='r9.6' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>.
<LineBreak/><LineBreak/>
Determine the specific unit type for any units identified per 
<InlineUIContainer><Button Content='r17.4' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>.
<LineBreak/><LineBreak/>
"@

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the e038 Orders Phase instructions in cell B39.
$ws.Range("B39").Value = $newText

# Move/select the active cell onto B39 (it was on A39 before the edit).
$ws.Range("B39").Select()
